$p = $ppt.ActivePresentation

# Slide 10: split the "Inherited condition..." paragraph into 4 paragraphs
$s10 = $p.Slides.Item(10)
$shape10 = $s10.Shapes.Item(2)
$tr10 = $shape10.TextFrame.TextRange
$tr10.Text = "Inherited condition in which an altered copy of the CDH1 gene is passed from generation to generation`rPresence of the gene can be detected by genetic testing`rAffected person can pass the gene to (on average) half of their children`rAffected persons carry the CDH1"

# Slide 11: split the "Affected individuals have microscopic cancers..." paragraph into 3 paragraphs
$s11 = $p.Slides.Item(11)
$shape11 = $s11.Shapes.Item(2)
$tr11 = $shape11.TextFrame.TextRange
$tr11.Text = "Affected individuals have microscopic cancers beginning to form in the top layer of the stomach at an early age`rMajority of affected individuals will develop visible cancer by age 40`rBy age 80, 70% of men and 56-83% of women are estimated to be at risk to develop visible cancer however some recent studies place this risk at 50%/33%"

# Slide 12: split the "CDH1 carriers with visible cancer..." paragraph into 2 paragraphs
# (note: the second sentence gains " at the time of diagnosis" at the end)
$s12 = $p.Slides.Item(12)
$shape12 = $s12.Shapes.Item(2)
$tr12 = $shape12.TextFrame.TextRange
$tr12.Text = "CDH1 carriers with visible cancer are termed “clinically apparent”`rCDH1 carriers with clinically apparent cancers which are large enough to cause symptoms generally are likely to have spread to lymph nodes at the time of diagnosis"
